$d = $word.ActiveDocument
# --- 1) A Gero 5559393 : add lastRenderedPageBreak ---
$r = $d.Content
$null = $r.Find.Execute("A Gero 5559393", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$xmlGero = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="01D508E1" w14:textId="77777777" w:rsidR="00772C40" w:rsidRPr="00772C40" w:rsidRDefault="00772C40" w:rsidP="002C4708"><w:pPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr></w:pPr><w:r w:rsidRPr="00772C40"><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:noProof/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>A Gero 5559393</w:t></w:r></w:p>
'@
$r.Paragraphs(1).Range.InsertXML($xmlGero)

# --- 2) A Simo 02/987665544 : remove lastRenderedPageBreak ---
$r = $d.Content
$null = $r.Find.Execute("A Simo 02/987665544", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$xmlSimo = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="64C24A25" w14:textId="77777777" w:rsidR="00772C40" w:rsidRPr="00772C40" w:rsidRDefault="00772C40" w:rsidP="002C4708"><w:pPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr></w:pPr><w:r w:rsidRPr="00772C40"><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:noProof/><w:lang w:val="en-US"/></w:rPr><w:t>A Simo 02/987665544</w:t></w:r></w:p>
'@
$r.Paragraphs(1).Range.InsertXML($xmlSimo)

# --- 3) Contact RoYaL does not exist. : add lastRenderedPageBreak ---
$r = $d.Content
$null = $r.Find.Execute("Contact RoYaL does not exist.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$xmlContactRoyal = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7CD202A1" w14:textId="77777777" w:rsidR="00772C40" w:rsidRPr="00772C40" w:rsidRDefault="00772C40" w:rsidP="002C4708"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr></w:pPr><w:r w:rsidRPr="00772C40"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:noProof/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Contact RoYaL does not exist.</w:t></w:r></w:p>
'@
$r.Paragraphs(1).Range.InsertXML($xmlContactRoyal)

# --- 4) RoYaL(Ivan) -> 666 : remove lastRenderedPageBreak ---
$r = $d.Content
$null = $r.Find.Execute("RoYaL(Ivan) -> 666", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$xmlRoyalIvan666 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1FBBDE37" w14:textId="77777777" w:rsidR="00772C40" w:rsidRPr="00772C40" w:rsidRDefault="00772C40" w:rsidP="002C4708"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:bCs/><w:noProof/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00772C40"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:noProof/><w:lang w:val="en-US"/></w:rPr><w:t>RoYaL(Ivan) -&gt; 666</w:t></w:r></w:p>
'@
$r.Paragraphs(1).Range.InsertXML($xmlRoyalIvan666)

# --- 5) Важно : add lastRenderedPageBreak (first run of paragraph) ---
$r = $d.Content
$null = $r.Find.Execute("Важно", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$xmlVazhno = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5AA80B3A" w14:textId="77777777" w:rsidR="00772C40" w:rsidRPr="00772C40" w:rsidRDefault="00772C40" w:rsidP="00772C40"><w:pPr><w:rPr><w:i/></w:rPr></w:pPr><w:r w:rsidRPr="00772C40"><w:rPr><w:i/></w:rPr><w:lastRenderedPageBreak/><w:t>Важно</w:t></w:r><w:r w:rsidRPr="00772C40"><w:rPr><w:i/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidRPr="00772C40"><w:rPr><w:i/></w:rPr><w:t>Общата цена се изчислява на базата на най</w:t></w:r><w:r w:rsidRPr="00772C40"><w:rPr><w:i/><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r><w:r w:rsidRPr="00772C40"><w:rPr><w:i/></w:rPr><w:t>новата цена за всеки продукт</w:t></w:r></w:p>
'@
$r.Paragraphs(1).Range.InsertXML($xmlVazhno)

# --- 6) Вход (heading after Важно paragraph) : remove lastRenderedPageBreak ---
$r = $d.Content
$null = $r.Find.Execute("Важно", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pVazhno = $r.Paragraphs(1)
$pVhod = $pVazhno.Next()
$xmlVhod = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7BC3912D" w14:textId="77777777" w:rsidR="00772C40" w:rsidRPr="00772C40" w:rsidRDefault="00772C40" w:rsidP="00772C40"><w:pPr><w:pStyle w:val="3"/><w:rPr><w:lang w:val="bg-BG"/></w:rPr></w:pPr><w:r w:rsidRPr="00772C40"><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>Вход</w:t></w:r></w:p>
'@
$pVhod.Range.InsertXML($xmlVhod)

